$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.25660000000001
$ws.Range("B4").Value = 4.748400000000006

$ws.Range("B5").Value = 5.305599999999998

$ws.Range("A6").Value = -21.58320000000003
$ws.Range("B6").Value = 5.562599999999999

$ws.Range("A7").Value = -21.38200000000001

$ws.Range("A8").Value = -21.56060000000002
$ws.Range("B8").Value = 5.058900000000001

$ws.Range("A16").Value = -21.62470000000003
$ws.Range("B16").Value = 5.194300000000005

$ws.Range("A20").Value = -22.09490000000002

$ws.Range("A21").Value = -20.48829999999999

$ws.Range("B22").Value = 5.782500000000002
